$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix typo: "ptions" -> "options" in heading text and its bookmark name.
# ---------------------------------------------------------------------------
$bmOld = $d.Bookmarks("some-general-ptions-not-discussed")
$bmRange = $bmOld.Range
$bmRange.Text = "Some general options not discussed"
$bmRangeDup = $bmOld.Range.Duplicate
$bmOld.Delete()
$d.Bookmarks.Add("some-general-options-not-discussed", $bmRangeDup) | Out-Null

# ---------------------------------------------------------------------------
# Helper functions used while appending the new material at the end of the
# document body (just before the final section break).
# ---------------------------------------------------------------------------

function New-Para($style) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $r.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Style = $style
    # Clear the auto-created empty run's text so we can build it up run by run.
    $p.Range.Text = ""
    return $p
}

function Add-Run($para, $text, $verbatim) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $r.InsertAfter($text)
    if ($verbatim) {
        $newEnd = $d.Content.End
        $runRange = $d.Range($newEnd - $text.Length, $newEnd)
        $runRange.Style = "VerbatimChar"
    }
}

# ---------------------------------------------------------------------------
# 2) New "BodyText" paragraph describing multiple output formats.
# ---------------------------------------------------------------------------
$p = New-Para "BodyText"
Add-Run $p "e.g. adding/updating" $false
Add-Run $p " " $false
Add-Run $p "output:   pdf_document: default   html_notebook: default   html_document:     df_print: paged   word_document: default" $true
Add-Run $p " " $false
Add-Run $p "in the YAML header will produce a a pdf document, and html notebook, and an html document, and a Microsoft Word document, all with one click of the Knit button!" $false

# ---------------------------------------------------------------------------
# 3) Bulleted list describing each output type. The list re-uses the same
#    bullet definition (abstractNum 991) already used elsewhere ("Compact"
#    style numbered/bulleted paragraphs).
# ---------------------------------------------------------------------------
$sample = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Caching:*") {
        $sample = $para
        break
    }
}
$baseTemplate = $sample.Range.ListFormat.ListTemplate

# -- pdf_document (starts a new list -> new numId) --------------------------
$p = New-Para "Compact"
Add-Run $p "pdf_document" $true
Add-Run $p ": Excellent for distributing a static (unchanging) document. Also great for submitting homework solutions to your instructor!" $false
$p.Range.ListFormat.ApplyListTemplateWithLevel($baseTemplate, $false, 1, $false, 0)
$topTemplate = $p.Range.ListFormat.ListTemplate

# -- html_document (continues the list started above) -----------------------
$p = New-Para "Compact"
Add-Run $p "html_document" $true
Add-Run $p ": Excellent for sharing interactive files. Can be viewed in a web browser." $false
$p.Range.ListFormat.ApplyListTemplateWithLevel($topTemplate, $true, 1, $false, 0)

# -- html_notebook (continues the list started above) ------------------------
$p = New-Para "Compact"
Add-Run $p "html_notebook" $true
Add-Run $p ": The excellent" $false
Add-Run $p " " $false
Add-Run $p [char]0x201C $false
Add-Run $p "R notebook" $false
Add-Run $p [char]0x201D $false
Add-Run $p "." $false
$p.Range.ListFormat.ApplyListTemplateWithLevel($topTemplate, $true, 1, $false, 0)

# -- nested bullet 1 (new sub-list -> new numId, level 2) --------------------
$p = New-Para "Compact"
Add-Run $p ("It" + [char]0x2019 + "s like an") $false
Add-Run $p " " $false
Add-Run $p "html_document" $true
Add-Run $p " " $false
Add-Run $p "but can be edited directly using R Studio and it can be used to reproduce the source code of the file generating it." $false
$p.Range.ListFormat.ApplyListTemplateWithLevel($baseTemplate, $false, 1, $false, 0)
$p.Range.ListFormat.ListLevelNumber = 2
$subTemplate = $p.Range.ListFormat.ListTemplate

# -- nested bullet 2 (continues the sub-list started above) ------------------
$p = New-Para "Compact"
Add-Run $p "Great to share with collaborators because they can update the document, add analysis, etc." $false
$p.Range.ListFormat.ApplyListTemplateWithLevel($subTemplate, $true, 1, $false, 0)
$p.Range.ListFormat.ListLevelNumber = 2

# -- word_document (back to the top-level list) ------------------------------
$p = New-Para "Compact"
Add-Run $p "word_document" $true
Add-Run $p ": If you have to have Microsoft Word editing capabilities, for some reason." $false
$p.Range.ListFormat.ApplyListTemplateWithLevel($topTemplate, $true, 1, $false, 0)

# ---------------------------------------------------------------------------
# 4) "Your turn" Heading3 with bookmark "your-turn-4".
# ---------------------------------------------------------------------------
$p = New-Para "Heading3"
Add-Run $p "Your turn" $false
$headingStart = $p.Range.Start
$headingEnd = $p.Range.End
# Trim the paragraph mark off the end of the range used for the bookmark.
$bmRange2 = $d.Range($headingStart, $headingEnd - 1)
$d.Bookmarks.Add("your-turn-4", $bmRange2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "FirstParagraph" instructions paragraph.
# ---------------------------------------------------------------------------
$p = New-Para "FirstParagraph"
Add-Run $p "Execute the command" $false
Add-Run $p " " $false
Add-Run $p 'download.file("https://github.com/jfrench/DataWrangleViz/blob/master/02-crash-course-in-rmd.nb.html")' $true
Add-Run $p " " $false
Add-Run $p "in the Console." $false

# ---------------------------------------------------------------------------
# 6) Two closing "BodyText" paragraphs.
# ---------------------------------------------------------------------------
$p = New-Para "BodyText"
Add-Run $p "Close R Studio." $false

$p = New-Para "BodyText"
Add-Run $p "Find the downloaded file and double-click to open it in R Studio." $false

Write-Output "Edit complete"
